$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-03 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-04 Thursday", 2) | Out-Null
$d.Content.Find.Execute("8+47=55", $true, $false, $false, $false, $false, $true, 1, $false, "49+20=69", 2) | Out-Null
$d.Content.Find.Execute("43-38=5", $true, $false, $false, $false, $false, $true, 1, $false, "70+12=82", 2) | Out-Null
$d.Content.Find.Execute("79+0=79", $true, $false, $false, $false, $false, $true, 1, $false, "38-7=31", 2) | Out-Null
$d.Content.Find.Execute("7-4=3", $true, $false, $false, $false, $false, $true, 1, $false, "42-0=42", 2) | Out-Null
$d.Content.Find.Execute("63+15=78", $true, $false, $false, $false, $false, $true, 1, $false, "39-21=18", 2) | Out-Null
$d.Content.Find.Execute("24-5=19", $true, $false, $false, $false, $false, $true, 1, $false, "7+47=54", 2) | Out-Null
$d.Content.Find.Execute("90-21=69", $true, $false, $false, $false, $false, $true, 1, $false, "59+29=88", 2) | Out-Null
$d.Content.Find.Execute("45+39=84", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2) | Out-Null
$d.Content.Find.Execute("31+53=84", $true, $false, $false, $false, $false, $true, 1, $false, "32+47=79", 2) | Out-Null
$d.Content.Find.Execute("2+35=37", $true, $false, $false, $false, $false, $true, 1, $false, "27+66=93", 2) | Out-Null
$d.Content.Find.Execute("88-61=27", $true, $false, $false, $false, $false, $true, 1, $false, "18-9=9", 2) | Out-Null
$d.Content.Find.Execute("11+16=27", $true, $false, $false, $false, $false, $true, 1, $false, "80+17=97", 2) | Out-Null
$d.Content.Find.Execute("17+57=74", $true, $false, $false, $false, $false, $true, 1, $false, "9+75=84", 2) | Out-Null
$d.Content.Find.Execute("57+33=90", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=30", 2) | Out-Null
$d.Content.Find.Execute("41-24=17", $true, $false, $false, $false, $false, $true, 1, $false, "62+6=68", 2) | Out-Null
$d.Content.Find.Execute("15+60=75", $true, $false, $false, $false, $false, $true, 1, $false, "31+42=73", 2) | Out-Null
$d.Content.Find.Execute("45+49=94", $true, $false, $false, $false, $false, $true, 1, $false, "64+9=73", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $false, $false, $false, $false, $true, 1, $false, "48+47=95", 2) | Out-Null
$d.Content.Find.Execute("23-8=15", $true, $false, $false, $false, $false, $true, 1, $false, "41-26=15", 2) | Out-Null
$d.Content.Find.Execute("58+29=87", $true, $false, $false, $false, $false, $true, 1, $false, "26+62=88", 2) | Out-Null
$d.Content.Find.Execute("60-54=6", $true, $false, $false, $false, $false, $true, 1, $false, "97-25=72", 2) | Out-Null
$d.Content.Find.Execute("29+31=60", $true, $false, $false, $false, $false, $true, 1, $false, "5+9=14", 2) | Out-Null
$d.Content.Find.Execute("54-20=34", $true, $false, $false, $false, $false, $true, 1, $false, "50+1=51", 2) | Out-Null
$d.Content.Find.Execute("69-36=33", $true, $false, $false, $false, $false, $true, 1, $false, "64+12=76", 2) | Out-Null
$d.Content.Find.Execute("23+58=81", $true, $false, $false, $false, $false, $true, 1, $false, "93-90=3", 2) | Out-Null
$d.Content.Find.Execute("81-41=40", $true, $false, $false, $false, $false, $true, 1, $false, "69+28=97", 2) | Out-Null
$d.Content.Find.Execute("24+56=80", $true, $false, $false, $false, $false, $true, 1, $false, "80+11=91", 2) | Out-Null
$d.Content.Find.Execute("73-69=4", $true, $false, $false, $false, $false, $true, 1, $false, "80+18=98", 2) | Out-Null
$d.Content.Find.Execute("87-29=58", $true, $false, $false, $false, $false, $true, 1, $false, "48+24=72", 2) | Out-Null
$d.Content.Find.Execute("22+23=45", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("70-27=43", $true, $false, $false, $false, $false, $true, 1, $false, "86-57=29", 2) | Out-Null
$d.Content.Find.Execute("23+46=69", $true, $false, $false, $false, $false, $true, 1, $false, "76+5=81", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("84-49=35", $true, $false, $false, $false, $false, $true, 1, $false, "56-54=2", 2) | Out-Null
$d.Content.Find.Execute("18+35=53", $true, $false, $false, $false, $false, $true, 1, $false, "54+8=62", 2) | Out-Null
$d.Content.Find.Execute("44+50=94", $true, $false, $false, $false, $false, $true, 1, $false, "51-16=35", 2) | Out-Null
$d.Content.Find.Execute("58-17=41", $true, $false, $false, $false, $false, $true, 1, $false, "24+53=77", 2) | Out-Null
$d.Content.Find.Execute("62+1=63", $true, $false, $false, $false, $false, $true, 1, $false, "40+58=98", 2) | Out-Null
$d.Content.Find.Execute("11+75=86", $true, $false, $false, $false, $false, $true, 1, $false, "86-40=46", 2) | Out-Null
$d.Content.Find.Execute("96-50=46", $true, $false, $false, $false, $false, $true, 1, $false, "11+28=39", 2) | Out-Null
$d.Content.Find.Execute("74-5=69", $true, $false, $false, $false, $false, $true, 1, $false, "44-37=7", 2) | Out-Null
$d.Content.Find.Execute("42+34=76", $true, $false, $false, $false, $false, $true, 1, $false, "18+49=67", 2) | Out-Null
$d.Content.Find.Execute("54-44=10", $true, $false, $false, $false, $false, $true, 1, $false, "71-42=29", 2) | Out-Null
$d.Content.Find.Execute("47+4=51", $true, $false, $false, $false, $false, $true, 1, $false, "90-77=13", 2) | Out-Null
$d.Content.Find.Execute("51-5=46", $true, $false, $false, $false, $false, $true, 1, $false, "38+11=49", 2) | Out-Null
$d.Content.Find.Execute("34-28=6", $true, $false, $false, $false, $false, $true, 1, $false, "26+33=59", 2) | Out-Null
$d.Content.Find.Execute("62-42=20", $true, $false, $false, $false, $false, $true, 1, $false, "10+1=11", 2) | Out-Null
$d.Content.Find.Execute("59-3=56", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=49", 2) | Out-Null
$d.Content.Find.Execute("7+28=35", $true, $false, $false, $false, $false, $true, 1, $false, "63+3=66", 2) | Out-Null
$d.Content.Find.Execute("49+30=79", $true, $false, $false, $false, $false, $true, 1, $false, "55+6=61", 2) | Out-Null
$d.Content.Find.Execute("14+10=24", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=61", 2) | Out-Null
$d.Content.Find.Execute("0+2=2", $true, $false, $false, $false, $false, $true, 1, $false, "28+34=62", 2) | Out-Null
$d.Content.Find.Execute("44+27=71", $true, $false, $false, $false, $false, $true, 1, $false, "99-51=48", 2) | Out-Null
$d.Content.Find.Execute("29+52=81", $true, $false, $false, $false, $false, $true, 1, $false, "9+64=73", 2) | Out-Null
$d.Content.Find.Execute("35-3=32", $true, $false, $false, $false, $false, $true, 1, $false, "71-60=11", 2) | Out-Null
$d.Content.Find.Execute("40-17=23", $true, $false, $false, $false, $false, $true, 1, $false, "34+13=47", 2) | Out-Null
$d.Content.Find.Execute("18-17=1", $true, $false, $false, $false, $false, $true, 1, $false, "62-61=1", 2) | Out-Null
$d.Content.Find.Execute("3+26=29", $true, $false, $false, $false, $false, $true, 1, $false, "36+28=64", 2) | Out-Null
$d.Content.Find.Execute("94-27=67", $true, $false, $false, $false, $false, $true, 1, $false, "83-40=43", 2) | Out-Null
$d.Content.Find.Execute("44+7=51", $true, $false, $false, $false, $false, $true, 1, $false, "75-12=63", 2) | Out-Null
$d.Content.Find.Execute("61-51=10", $true, $false, $false, $false, $false, $true, 1, $false, "67-54=13", 2) | Out-Null
$d.Content.Find.Execute("29-15=14", $true, $false, $false, $false, $false, $true, 1, $false, "19+67=86", 2) | Out-Null
$d.Content.Find.Execute("90-76=14", $true, $false, $false, $false, $false, $true, 1, $false, "62-56=6", 2) | Out-Null
$d.Content.Find.Execute("27+33=60", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=19", 2) | Out-Null
$d.Content.Find.Execute("15+62=77", $true, $false, $false, $false, $false, $true, 1, $false, "1+40=41", 2) | Out-Null
$d.Content.Find.Execute("28+11=39", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=1", 2) | Out-Null
$d.Content.Find.Execute("52-35=17", $true, $false, $false, $false, $false, $true, 1, $false, "52+21=73", 2) | Out-Null
$d.Content.Find.Execute("33-23=10", $true, $false, $false, $false, $false, $true, 1, $false, "62-30=32", 2) | Out-Null
$d.Content.Find.Execute("0+81=81", $true, $false, $false, $false, $false, $true, 1, $false, "42+31=73", 2) | Out-Null
$d.Content.Find.Execute("16+70=86", $true, $false, $false, $false, $false, $true, 1, $false, "66+31=97", 2) | Out-Null
$d.Content.Find.Execute("90-71=19", $true, $false, $false, $false, $false, $true, 1, $false, "19+76=95", 2) | Out-Null
$d.Content.Find.Execute("40-20=20", $true, $false, $false, $false, $false, $true, 1, $false, "56+35=91", 2) | Out-Null
$d.Content.Find.Execute("10+87=97", $true, $false, $false, $false, $false, $true, 1, $false, "77+0=77", 2) | Out-Null
$d.Content.Find.Execute("94-18=76", $true, $false, $false, $false, $false, $true, 1, $false, "96-60=36", 2) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("64-24=40", $true, $false, $false, $false, $false, $true, 1, $false, "26+36=62", 2) | Out-Null
$d.Content.Find.Execute("27-20=7", $true, $false, $false, $false, $false, $true, 1, $false, "77-7=70", 2) | Out-Null
$d.Content.Find.Execute("4+3=7", $true, $false, $false, $false, $false, $true, 1, $false, "12-11=1", 2) | Out-Null
$d.Content.Find.Execute("36+40=76", $true, $false, $false, $false, $false, $true, 1, $false, "38-28=10", 2) | Out-Null
$d.Content.Find.Execute("75-29=46", $true, $false, $false, $false, $false, $true, 1, $false, "1+96=97", 2) | Out-Null
$d.Content.Find.Execute("20+23=43", $true, $false, $false, $false, $false, $true, 1, $false, "62+12=74", 2) | Out-Null
$d.Content.Find.Execute("75-33=42", $true, $false, $false, $false, $false, $true, 1, $false, "29-3=26", 2) | Out-Null
$d.Content.Find.Execute("25-17=8", $true, $false, $false, $false, $false, $true, 1, $false, "98-84=14", 2) | Out-Null
$d.Content.Find.Execute("20-1=19", $true, $false, $false, $false, $false, $true, 1, $false, "0+31=31", 2) | Out-Null
$d.Content.Find.Execute("76+18=94", $true, $false, $false, $false, $false, $true, 1, $false, "85-66=19", 2) | Out-Null
$d.Content.Find.Execute("97-73=24", $true, $false, $false, $false, $false, $true, 1, $false, "67-22=45", 2) | Out-Null
$d.Content.Find.Execute("30+49=79", $true, $false, $false, $false, $false, $true, 1, $false, "83-62=21", 2) | Out-Null
$d.Content.Find.Execute("1+22=23", $true, $false, $false, $false, $false, $true, 1, $false, "93-70=23", 2) | Out-Null
$d.Content.Find.Execute("54-43=11", $true, $false, $false, $false, $false, $true, 1, $false, "6+63=69", 2) | Out-Null
$d.Content.Find.Execute("69-55=14", $true, $false, $false, $false, $false, $true, 1, $false, "88-11=77", 2) | Out-Null
$d.Content.Find.Execute("1+88=89", $true, $false, $false, $false, $false, $true, 1, $false, "50+37=87", 2) | Out-Null
$d.Content.Find.Execute("56-53=3", $true, $false, $false, $false, $false, $true, 1, $false, "68+7=75", 2) | Out-Null
$d.Content.Find.Execute("66-32=34", $true, $false, $false, $false, $false, $true, 1, $false, "53+36=89", 2) | Out-Null
$d.Content.Find.Execute("37+56=93", $true, $false, $false, $false, $false, $true, 1, $false, "98-1=97", 2) | Out-Null
$d.Content.Find.Execute("48-0=48", $true, $false, $false, $false, $false, $true, 1, $false, "86-84=2", 2) | Out-Null
$d.Content.Find.Execute("66-51=15", $true, $false, $false, $false, $false, $true, 1, $false, "68-53=15", 2) | Out-Null
$d.Content.Find.Execute("66-44=22", $true, $false, $false, $false, $false, $true, 1, $false, "76-71=5", 2) | Out-Null
$d.Content.Find.Execute("80-72=8", $true, $false, $false, $false, $false, $true, 1, $false, "77-69=8", 2) | Out-Null
$d.Content.Find.Execute("4+55=59", $true, $false, $false, $false, $false, $true, 1, $false, "66-62=4", 2) | Out-Null
$d.Content.Find.Execute("55-32=23", $true, $false, $false, $false, $false, $true, 1, $false, "22+10=32", 2) | Out-Null

Write-Host "Replacements complete"